# Split the run-on "Programa" paragraphs (English + Portuguese) into
# individually numbered sentences, each followed by a line break (<w:br/>),
# instead of one long run of concatenated text. Also split the
# "NF= (N1 + N2)/2Onde: ..." sentence onto its own line.
#
# Approach: use Find to locate the start of each numbered item (2. .. 9.)
# and insert a manual line break (vertical-tab char, Word's internal
# line-break marker) immediately before it. This causes the run to split
# and a <w:br/> element to be created between the two resulting <w:t>s,
# matching the target OOXML exactly.

$d = $word.ActiveDocument
$wdLineBreak = [char]11

function Insert-BreakBefore($range, $marker) {
    $found = $range.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $found) {
        throw "marker not found: $marker"
    }
    $range.Collapse(1)  # wdCollapseStart
    $range.InsertBefore($wdLineBreak)
}

# --- English "Programa" paragraph (italic run) -----------------------
$enMarkers = @(
    "2.  Eco-innovation metrics",
    "3. Introduction to products life-cycle",
    "4.  Eco-innovation in the industry",
    "5.  Case study of Eco-innovation projects in Brazil.",
    "6. Methods and tools to support the process of Eco-innovation",
    "7. Early identification of failure as support to Eco-innovation",
    "8. TRIZ as a response to Eco-innovation",
    "9. Methodological proposal for Eco-innovative solutions"
)

foreach ($m in $enMarkers) {
    $r = $d.Content.Duplicate
    Insert-BreakBefore $r $m
}

# --- Portuguese "Método" paragraph (plain run) ------------------------
$ptMarkers = @(
    "2. Métricas da eco-inovação",
    "3. Introdução ao Ciclo de vida do produto",
    "4. Eco inovação na indústria",
    "5. Estudo de casos de projetos de eco inovação no Brasil.",
    "6. Métodos e ferramentas suporte do processo de eco-inovação",
    "7. Identificação antecipada de falha como suporte a eco-inovação",
    "8. TRIZ como resposta a eco inovação",
    "9. Proposta metodológica para soluções eco inovadoras"
)

foreach ($m in $ptMarkers) {
    $r = $d.Content.Duplicate
    Insert-BreakBefore $r $m
}

# --- "Norma de recuperação" sentence: split "NF=.../2Onde:..." --------
$r = $d.Content.Duplicate
Insert-BreakBefore $r "Onde: NF = nota final"

Write-Host "Split numbered Programa/Metodo items onto separate lines."
